# Auto-generated edit script: updates cryptos list values to match the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.335.16'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = '2.058.26'
$ws.Range('E3').Value = '  +1.34%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.79'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('E6').Value = '  +2.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '57.83'
$ws.Range('E7').Value = '  +5.09%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.383'
$ws.Range('E9').Value = '  +3.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.31'
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('E12').Value = '  +1.39%  '
$ws.Range('D13').Value = '2.363.29'
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.41'
$ws.Range('E14').Value = '  +1.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.81'
$ws.Range('E15').Value = '  +2.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.773'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.16'
$ws.Range('E17').Value = '  +1.42%  '
$ws.Range('D18').Value = '2.053.50'
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('D19').Value = '37.284.24'
$ws.Range('E19').Value = '  +1.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.24'
$ws.Range('E20').Value = '  +14.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '69.22'
$ws.Range('E21').Value = '  +2.21%  '
$ws.Range('D22').Value = '0.0₃0811'
$ws.Range('E22').Value = '  +1.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '225.16'
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('E25').Value = '  +2.63%  '
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.50'
$ws.Range('E27').Value = '  +1.61%  '
$ws.Range('E28').Value = '  +6.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.80'
$ws.Range('E29').Value = '  +2.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.128'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.08'
$ws.Range('E31').Value = '  +0.51%  '
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('E33').Value = '  +2.65%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0616'
$ws.Range('E34').Value = '  +2.22%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.55'
$ws.Range('E35').Value = '  +3.19%  '
$ws.Range('E36').Value = '  +6.62%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.85'
$ws.Range('E38').Value = '  +1.20%  '
$ws.Range('B39').Value = 'WEMIXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.74'
$ws.Range('E39').Value = '  -1.44%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.26'
$ws.Range('E40').Value = '  -1.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.69'
$ws.Range('E41').Value = '  +12.88%  '
$ws.Range('E42').Value = '  +1.14%  '
$ws.Range('D43').Value = '1.478.47'
$ws.Range('E43').Value = '  +1.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '96.66'
$ws.Range('E44').Value = '  +6.53%  '
$ws.Range('E45').Value = '  -1.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.16'
$ws.Range('E46').Value = '  +4.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0210'
$ws.Range('E47').Value = '  +3.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.49'
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('E49').Value = '  +1.28%  '
$ws.Range('E50').Value = '  +4.10%  '
$ws.Range('E51').Value = '  +1.79%  '
